$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range('D2').Value = '72.322.71'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.641.70'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '583.54'
$ws.Range('E5').Value = '  -3.15%  '
$ws.Range('D6').Value = '174.66'
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.519'
$ws.Range('E8').Value = '  -0.87%  '
$ws.Range('D9').Value = '2.641.67'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('D11').Value = '0.170'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '4.91'
$ws.Range('E13').Value = '  -2.55%  '
$ws.Range('D14').Value = '3.123.83'
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0000185'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '72.032.06'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '25.85'
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '2.623.69'
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('D19').Value = '8.43'
$ws.Range('E19').Value = '  +5.69%  '
$ws.Range('D20').Value = '12.12'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').Value = '373.53'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').Value = '2.04'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '70.76'
$ws.Range('E25').Value = '  -2.25%  '
$ws.Range('D26').Value = '4.24'
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('E27').Value = '  -3.41%  '
$ws.Range('D28').Value = '2.776.91'
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = '0.0₃0948'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '7.95'
$ws.Range('D32').Value = '496.65'
$ws.Range('E32').Value = '  -4.17%  '
$ws.Range('E33').Value = '  -2.58%  '
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '162.74'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').Value = '19.19'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('E38').Value = '  +3.11%  '
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('D40').Value = '1.35'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('E42').Value = '  -6.67%  '
$ws.Range('D43').Value = '2.56'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').Value = '4.88'
$ws.Range('E44').Value = '  -2.98%  '
$ws.Range('D45').Value = '0.326'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('D46').Value = '39.01'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').Value = '151.93'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('E48').Value = '  -2.25%  '
$ws.Range('D49').Value = '0.544'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('D51').Value = '0.601'
$ws.Range('E51').Value = '  -0.42%  '

$textRange.Style = "Normal"
